$d = $word.ActiveDocument

function Get-ParagraphByExactText($doc, $exactText) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        # Paragraphs end with a paragraph mark (CR, chr 13) - trim it for comparison.
        $trimmed = $t.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $exactText) {
            return $p
        }
    }
    return $null
}

# The three consecutive Manager bullet items before the edit:
#   A: "добавлять авиакомпании и редактировать информацию о них;"
#   B: "просматривать рейсы и доступные билеты;"
#   C: "оформлять подажу билета для клиента (на стойке регистрации);"
#
# The edit folds B's wording into A (replacing A's old text) and removes the
# now-redundant B and C bullet paragraphs entirely.

# Delete paragraph C first, then B, so earlier paragraph indices/objects stay valid.
$pC = Get-ParagraphByExactText $d "оформлять подажу билета для клиента (на стойке регистрации);"
if ($pC -ne $null) {
    $pC.Range.Delete()
}

$pB = Get-ParagraphByExactText $d "просматривать рейсы и доступные билеты;"
if ($pB -ne $null) {
    $pB.Range.Delete()
}

# Replace paragraph A's run text with B's old wording.
$pA = Get-ParagraphByExactText $d "добавлять авиакомпании и редактировать информацию о них;"
$pA.Range.Find.Execute("добавлять авиакомпании и редактировать информацию о них", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "просматривать рейсы и доступные билеты", 2)

# Re-locate paragraph A (its object may have been reseated by the edit above) and
# split "доступные билеты" into its own run, matching the target's run layout,
# by toggling a character property on and back off (a formatting no-op that
# still forces Word to materialize a separate run at that boundary).
$pA2 = Get-ParagraphByExactText $d "просматривать рейсы и доступные билеты;"
$sub = $pA2.Range
$sub.Find.Execute("доступные билеты", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sub.Bold = 1
$sub.Bold = 0
